$d = $word.ActiveDocument

# Locate the paragraph containing the "{m:userdoc 'zone1'}" template token
# (originally split as two runs: "{m" and ":userdoc 'zone1'}") and rewrite
# it as four runs: "{", "m", ":userdoc 'zone1'", "}".
#
# A plain Find/Replace (or Range.Text / InsertAfter) collapses the
# surrounding text back into a single run because the formatting is
# identical on both sides of the split point. To force Word to keep the
# runs separate without touching any character formatting, we briefly drop
# a bookmark at each split point (which always breaks the enclosing run in
# two) and then delete that bookmark again - the run boundary it created
# survives the bookmark's removal.
function Split-RunAt($pos) {
    $name = "m2docSplitMarker"
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $r) | Out-Null
    $d.Bookmarks($name).Delete()
}

$target = "{m:userdoc 'zone1'}"

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($target)) {
        $para = $p
        break
    }
}

$start = $para.Range.Start

# Split points (character offsets from the start of the document):
#   "{" | "m" | ":userdoc 'zone1'" | "}"
$splitAfterBrace = $start + 1   # after "{"
$splitAfterM     = $start + 2   # after "m"
$splitAfterQuote = $start + 18  # after ":userdoc 'zone1'"

# Apply splits from right to left so earlier offsets stay valid.
Split-RunAt $splitAfterQuote
Split-RunAt $splitAfterM
Split-RunAt $splitAfterBrace
